$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update working days (E5), leave-with-permission days (F5), violation count (G5)
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1

# Row 12: update approved leader score (L12) and total task count (N12)
$ws.Range("L12").Value = 80.0
$ws.Range("N12").Value = 5
